# Change the table style used by the three data tables in this deck from
# the built-in "No Style, No Grid" style to the built-in "No Style, Table
# Grid" style ({D32A397A-1999-4E44-B5A6-C4F98FDC5BF0} ->
# {57916217-4300-40F5-8B69-F576111CB368}).

$OldStyleId = "{D32A397A-1999-4E44-B5A6-C4F98FDC5BF0}"
$NewStyleId = "{57916217-4300-40F5-8B69-F576111CB368}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style -eq $OldStyleId) {
                # Table styles are read-only as a plain property assignment
                # in the PowerPoint object model - use ApplyStyle to change
                # which built-in/table style gallery entry is applied.
                $table.ApplyStyle($NewStyleId)
            }
        }
    }
}
